$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (H1:J1) -- "Start_Story_Id", "After_Story_Id", "Letter_Id"
$ws.Range("H1").Value = "Start_Story_Id"
$ws.Range("I1").Value = "After_Story_Id"
$ws.Range("J1").Value = "Letter_Id"
$ws.Range("H1:J1").Style = $ws.Range("G1").Style

# New data cells (H2:J3) -- "Test" placeholder values
$ws.Range("H2").Value = "Test"
$ws.Range("I2").Value = "Test"
$ws.Range("J2").Value = "Test"
$ws.Range("H3").Value = "Test"
$ws.Range("I3").Value = "Test"
$ws.Range("J3").Value = "Test"

# Give the new data cells a thin box border (vertical-center alignment only)
$rng = $ws.Range("H2:J3")
$rng.Borders.Color = 0
$rng.Borders.Weight = 2
$rng.Borders.LineStyle = 1

# Column J width
$ws.Columns.Item(10).ColumnWidth = 10.7

# Dimension / selection bookkeeping
$ws.Range("J10").Select()
